$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2554.0454
$ws.Range("I98").Value = 1925.25
$ws.Range("J98").Value = 4230.8335
$ws.Range("K98").Value = 1925.25
$ws.Range("L98").Value = 4230.8335
$ws.Range("M98").Value = -427.25
$ws.Range("N98").Value = -7226.8335

$ws.Range("H122").Value = 2554.0454
$ws.Range("I122").Value = 1925.25
$ws.Range("J122").Value = 4230.8335
$ws.Range("K122").Value = 5775.75
$ws.Range("L122").Value = 12692.5005
$ws.Range("M122").Value = -3325.75
$ws.Range("N122").Value = -17592.5005

$ws.Range("H137").Value = 1045.4546
$ws.Range("I137").Value = 1014.7407
$ws.Range("J137").Value = 1183.6666
$ws.Range("K137").Value = 3044.2221
$ws.Range("L137").Value = 3550.9998
$ws.Range("M137").Value = -494.2221
$ws.Range("N137").Value = -8650.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4830.517
$ws.Range("I32").Value = 5119.423
$ws.Range("J32").Value = 2326.6667
$ws.Range("K32").Value = 5119.423
$ws.Range("L32").Value = 2326.6667
$ws.Range("M32").Value = -4832.423
$ws.Range("N32").Value = -2900.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 103909
$ws.Range("J42").Value = 103909
$ws.Range("L42").Value = 103909
$ws.Range("N42").Value = -104565

$ws.Range("H134").Value = 3947.2856
$ws.Range("I134").Value = 1037.2903
$ws.Range("J134").Value = 26499.75
$ws.Range("K134").Value = 3111.8709
$ws.Range("L134").Value = 79499.25
$ws.Range("M134").Value = -576.8708999999999
$ws.Range("N134").Value = -84569.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 27028450
$ws.Range("J131").Value = 1790.1072
$ws.Range("L131").Value = 5370.321599999999
$ws.Range("N131").Value = -15450.3216

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 20000
$ws.Range("I52").Value = 15000
$ws.Range("J52").Value = 25000
$ws.Range("K52").Value = 15000
$ws.Range("L52").Value = 25000
$ws.Range("M52").Value = -14741
$ws.Range("N52").Value = -25518

$ws.Range("H58").Value = 540
$ws.Range("I58").Value = 540
$ws.Range("K58").Value = 540
$ws.Range("M58").Value = -263

$ws.Range("H80").Value = 3366.6667
$ws.Range("J80").Value = 5325
$ws.Range("L80").Value = 5325
$ws.Range("N80").Value = -7321

$ws.Range("H83").Value = 3366.6667
$ws.Range("J83").Value = 5325
$ws.Range("L83").Value = 26625
$ws.Range("N83").Value = -36609

$ws.Range("H107").Value = 816.5
$ws.Range("I107").Value = 869.3333
$ws.Range("J107").Value = 748.5714
$ws.Range("K107").Value = 869.3333
$ws.Range("L107").Value = 748.5714
$ws.Range("M107").Value = 1050.6667
$ws.Range("N107").Value = -4588.5714

$ws.Range("H132").Value = 3072.6875
$ws.Range("I132").Value = 2430.5833
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 7291.749899999999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4761.749899999999
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1199.7693
$ws.Range("I16").Value = 1091.9
$ws.Range("J16").Value = 1559.3334
$ws.Range("K16").Value = 1091.9
$ws.Range("L16").Value = 1559.3334
$ws.Range("M16").Value = -921.9000000000001
$ws.Range("N16").Value = -1899.3334

$ws.Range("H22").Value = 1666.5
$ws.Range("I22").Value = 1500
$ws.Range("J22").Value = 1999.5
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 1999.5
$ws.Range("M22").Value = -1205
$ws.Range("N22").Value = -2589.5

$ws.Range("H27").Value = 1666.5
$ws.Range("I27").Value = 1500
$ws.Range("J27").Value = 1999.5
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 1999.5
$ws.Range("M27").Value = -1393
$ws.Range("N27").Value = -2213.5

$ws.Range("H46").Value = 2300.4
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 2751
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 2751
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -3127

$ws.Range("H55").Value = 204.25
$ws.Range("I55").Value = 76.5
$ws.Range("J55").Value = 587.5
$ws.Range("K55").Value = 76.5
$ws.Range("L55").Value = 587.5
$ws.Range("M55").Value = 96.5
$ws.Range("N55").Value = -933.5

$ws.Range("H68").Value = 1780.9333
$ws.Range("I68").Value = 1709.5
$ws.Range("J68").Value = 2066.6667
$ws.Range("K68").Value = 1709.5
$ws.Range("L68").Value = 2066.6667
$ws.Range("M68").Value = -960.5
$ws.Range("N68").Value = -3564.6667

$ws.Range("H71").Value = 1780.9333
$ws.Range("I71").Value = 1709.5
$ws.Range("J71").Value = 2066.6667
$ws.Range("K71").Value = 8547.5
$ws.Range("L71").Value = 10333.3335
$ws.Range("M71").Value = -4803.5
$ws.Range("N71").Value = -17821.3335

$ws.Range("H82").Value = 795
$ws.Range("I82").Value = 792
$ws.Range("J82").Value = 796.5
$ws.Range("K82").Value = 792
$ws.Range("L82").Value = 796.5
$ws.Range("M82").Value = -431
$ws.Range("N82").Value = -1518.5

$ws.Range("H85").Value = 795
$ws.Range("I85").Value = 792
$ws.Range("J85").Value = 796.5
$ws.Range("K85").Value = 792
$ws.Range("L85").Value = 796.5
$ws.Range("M85").Value = 456
$ws.Range("N85").Value = -3292.5

$ws.Range("H109").Value = 24000
$ws.Range("J109").Value = 24000
$ws.Range("L109").Value = 24000
$ws.Range("N109").Value = -26774

$ws.Range("H132").Value = 79911.69500000001
$ws.Range("I132").Value = 3712.25
$ws.Range("J132").Value = 113778.11
$ws.Range("K132").Value = 11136.75
$ws.Range("L132").Value = 341334.33
$ws.Range("M132").Value = -8606.75
$ws.Range("N132").Value = -346394.33

$ws.Range("H133").Value = 44950
$ws.Range("J133").Value = 44950
$ws.Range("L133").Value = 44950
$ws.Range("N133").Value = -50010

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 41673800
$ws.Range("I62").Value = 45460692
$ws.Range("K62").Value = 45460692
$ws.Range("M62").Value = -45460068

$ws.Range("H65").Value = 41673800
$ws.Range("I65").Value = 45460692
$ws.Range("K65").Value = 227303460
$ws.Range("M65").Value = -227300340

$ws.Range("H132").Value = 2791.2122
$ws.Range("I132").Value = 2122.3704
$ws.Range("K132").Value = 6367.111199999999
$ws.Range("M132").Value = -3837.111199999999
